$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Halve the yearly base progress value, and bump the MOD2 multiplier ---
$ws.Range("P1").Value = 0.75
$ws.Range("L2").Value = 6

# --- Add a new "STD + MOD1" combined yearly-rate block (rows 6-9, cols K-N) ---
$ws.Range("K6").Value = "Yearly 2"
$ws.Range("L6").Value = 10
$ws.Range("M6").Value = "MOD1"
$ws.Range("N6").Value = 0.5

$ws.Range("M7").Value = "Yearly"
$ws.Range("N7").Formula = "=N6*L6"

$ws.Range("K9").Value = "Sum"
$ws.Range("M9").Value = "STD + MOD1"
$ws.Range("N9").Formula = "=N7+L2"

# --- Re-point the later table rows (native reform progress, rows 12-17) to use
#     the new combined "Sum" divisor ($N$9) instead of the old $L$2 / $N$4 ---
$ws.Range("E12").Formula = "=D12/`$N`$9"
$ws.Range("E13").Formula = "=D13/`$N`$9"
$ws.Range("E14").Formula = "=D14/`$N`$9"
$ws.Range("E15").Formula = "=D15/`$N`$9"
$ws.Range("E16").Formula = "=D16/`$N`$9"
$ws.Range("E17").Formula = "=D17/`$N`$9"

$ws.Range("G13").Formula = "=`$D13/`$N`$9"
$ws.Range("G14").Formula = "=`$D14/`$N`$9"
$ws.Range("G15").Formula = "=`$D15/`$N`$9"
$ws.Range("G16").Formula = "=`$D16/`$N`$9"
$ws.Range("G17").Formula = "=`$D17/`$N`$9"

# --- Cosmetic: widen column M a bit, and move the active selection ---
$ws.Range("M1").EntireColumn.ColumnWidth = 11.5
$ws.Range("O11").Select()
